# StyleTemplate.xlsx workbook update
#
# The original workbook has two sheets:
#   1) "StyleTemplate"  - a demo sheet driven by conditional formatting / dxfs
#   2) "StyleTemplate1" - the real style template content
#
# The commit collapses the workbook down to a single sheet: the old
# "StyleTemplate" sheet (and the shared string / conditional formatting it
# alone used) is removed, and the remaining "StyleTemplate1" sheet is
# renamed back to "StyleTemplate" (taking over as the first/only sheet).
# The selection on that sheet is also updated to F7.

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete this sheet" interactive prompt.
$excel.DisplayAlerts = $false

# Drop the first sheet ("StyleTemplate") entirely - this also drops the
# conditional formatting (and the shared string) that only it referenced.
$oldSheet = $wb.Worksheets.Item("StyleTemplate")
$oldSheet.Delete()

# The remaining sheet ("StyleTemplate1") becomes the new, sole
# "StyleTemplate" sheet.
$templateSheet = $wb.Worksheets.Item("StyleTemplate1")
$templateSheet.Name = "StyleTemplate"

# Make it the active sheet and move the selection to F7, matching the
# saved view state in the workbook.
$templateSheet.Activate()
$templateSheet.Range("F7").Select()

$excel.DisplayAlerts = $true
